$wb = $excel.ActiveWorkbook

# Sheet 1: Posterior mean
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = -0.49140043533939
$ws.Range("D2").Value = -0.7997642541265038
$ws.Range("E2").Value = -0.8206483095871268
$ws.Range("F2").Value = -0.3442825030917555
$ws.Range("G2").Value = -0.6755872914676799
$ws.Range("H2").Value = -0.8061194045143951
$ws.Range("I2").Value = -0.4622597767012113
$ws.Range("J2").Value = -0.1958741043508663
$ws.Range("K2").Value = -0.2727279655791058
$ws.Range("L2").Value = 0.200459520130608
$ws.Range("M2").Value = -0.7365968399799289
$ws.Range("B3").Value = -0.49140043533939
$ws.Range("D3").Value = 0.633631974797373
$ws.Range("E3").Value = 0.4909245962882581
$ws.Range("F3").Value = 0.5347378990608032
$ws.Range("G3").Value = 0.6966241224194505
$ws.Range("H3").Value = 0.6257441245420285
$ws.Range("I3").Value = 0.5978091670819844
$ws.Range("J3").Value = 0.367936686628241
$ws.Range("K3").Value = 0.1981121834566008
$ws.Range("L3").Value = 0.01638312656337109
$ws.Range("M3").Value = 0.6734294713048488
$ws.Range("B4").Value = -0.7997642541265038
$ws.Range("C4").Value = 0.633631974797373
$ws.Range("E4").Value = 0.8657763548452646
$ws.Range("F4").Value = 0.5853457934066713
$ws.Range("G4").Value = 0.8789655368041299
$ws.Range("H4").Value = 0.9975075667717772
$ws.Range("I4").Value = 0.6899436301868895
$ws.Range("J4").Value = 0.3775471292142287
$ws.Range("K4").Value = 0.2866524564229284
$ws.Range("L4").Value = -0.09393857860515666
$ws.Range("M4").Value = 0.9209543843751962
$ws.Range("B5").Value = -0.8206483095871268
$ws.Range("C5").Value = 0.4909245962882581
$ws.Range("D5").Value = 0.8657763548452646
$ws.Range("F5").Value = 0.2158915863641953
$ws.Range("G5").Value = 0.6737927938227366
$ws.Range("H5").Value = 0.8746474862674289
$ws.Range("I5").Value = 0.3818787843660253
$ws.Range("J5").Value = 0.1399592825917317
$ws.Range("K5").Value = 0.2869606147580508
$ws.Range("L5").Value = -0.221298315099816
$ws.Range("M5").Value = 0.7784781961892289
$ws.Range("B6").Value = -0.3442825030917555
$ws.Range("C6").Value = 0.5347378990608032
$ws.Range("D6").Value = 0.5853457934066713
$ws.Range("E6").Value = 0.2158915863641953
$ws.Range("G6").Value = 0.7463595442535067
$ws.Range("H6").Value = 0.5677451094252198
$ws.Range("I6").Value = 0.8857979565169948
$ws.Range("J6").Value = 0.5680329442337095
$ws.Range("K6").Value = 0.1376996001907283
$ws.Range("L6").Value = 0.1608548349944636
$ws.Range("M6").Value = 0.619177195407188
$ws.Range("B7").Value = -0.6755872914676799
$ws.Range("C7").Value = 0.6966241224194505
$ws.Range("D7").Value = 0.8789655368041299
$ws.Range("E7").Value = 0.6737927938227366
$ws.Range("F7").Value = 0.7463595442535067
$ws.Range("H7").Value = 0.8688983785236081
$ws.Range("I7").Value = 0.8343740705772256
$ws.Range("J7").Value = 0.507918638109366
$ws.Range("K7").Value = 0.2609492386575604
$ws.Range("L7").Value = 0.02148188021592801
$ws.Range("M7").Value = 0.9329719601552322
$ws.Range("B8").Value = -0.8061194045143951
$ws.Range("C8").Value = 0.6257441245420285
$ws.Range("D8").Value = 0.9975075667717772
$ws.Range("E8").Value = 0.8746474862674289
$ws.Range("F8").Value = 0.5677451094252198
$ws.Range("G8").Value = 0.8688983785236081
$ws.Range("I8").Value = 0.6746978651626455
$ws.Range("J8").Value = 0.3639545243475186
$ws.Range("K8").Value = 0.2881234356303738
$ws.Range("L8").Value = -0.1034214358551055
$ws.Range("M8").Value = 0.9160401918518573
$ws.Range("B9").Value = -0.4622597767012113
$ws.Range("C9").Value = 0.5978091670819844
$ws.Range("D9").Value = 0.6899436301868895
$ws.Range("E9").Value = 0.3818787843660253
$ws.Range("F9").Value = 0.8857979565169948
$ws.Range("G9").Value = 0.8343740705772256
$ws.Range("H9").Value = 0.6746978651626455
$ws.Range("J9").Value = 0.5604356004376135
$ws.Range("K9").Value = 0.1858880121523191
$ws.Range("L9").Value = 0.1215624509034098
$ws.Range("M9").Value = 0.7305792375421298
$ws.Range("B10").Value = -0.1958741043508663
$ws.Range("C10").Value = 0.367936686628241
$ws.Range("D10").Value = 0.3775471292142287
$ws.Range("E10").Value = 0.1399592825917317
$ws.Range("F10").Value = 0.5680329442337095
$ws.Range("G10").Value = 0.507918638109366
$ws.Range("H10").Value = 0.3639545243475186
$ws.Range("I10").Value = 0.5604356004376135
$ws.Range("K10").Value = 0.1095345970178467
$ws.Range("L10").Value = 0.1558609204711422
$ws.Range("M10").Value = 0.4292212589197415
$ws.Range("B11").Value = -0.2727279655791058
$ws.Range("C11").Value = 0.1981121834566008
$ws.Range("D11").Value = 0.2866524564229284
$ws.Range("E11").Value = 0.2869606147580508
$ws.Range("F11").Value = 0.1376996001907283
$ws.Range("G11").Value = 0.2609492386575604
$ws.Range("H11").Value = 0.2881234356303738
$ws.Range("I11").Value = 0.1858880121523191
$ws.Range("J11").Value = 0.1095345970178467
$ws.Range("L11").Value = -0.02389876597373879
$ws.Range("M11").Value = 0.2764023456118498
$ws.Range("B12").Value = 0.200459520130608
$ws.Range("C12").Value = 0.01638312656337109
$ws.Range("D12").Value = -0.09393857860515666
$ws.Range("E12").Value = -0.221298315099816
$ws.Range("F12").Value = 0.1608548349944636
$ws.Range("G12").Value = 0.02148188021592801
$ws.Range("H12").Value = -0.1034214358551055
$ws.Range("I12").Value = 0.1215624509034098
$ws.Range("J12").Value = 0.1558609204711422
$ws.Range("K12").Value = -0.02389876597373879
$ws.Range("M12").Value = -0.0397879275596625
$ws.Range("B13").Value = -0.7365968399799289
$ws.Range("C13").Value = 0.6734294713048488
$ws.Range("D13").Value = 0.9209543843751962
$ws.Range("E13").Value = 0.7784781961892289
$ws.Range("F13").Value = 0.619177195407188
$ws.Range("G13").Value = 0.9329719601552322
$ws.Range("H13").Value = 0.9160401918518573
$ws.Range("I13").Value = 0.7305792375421298
$ws.Range("J13").Value = 0.4292212589197415
$ws.Range("K13").Value = 0.2764023456118498
$ws.Range("L13").Value = -0.0397879275596625

# Sheet 2: Pr(x>0)
$ws = $wb.Worksheets.Item(2)
$ws.Range("C2").Value = 0.2374375
$ws.Range("D2").Value = 0.09193750000000001
$ws.Range("E2").Value = 0.04259375
$ws.Range("F2").Value = 0.19378125
$ws.Range("G2").Value = 0.14053125
$ws.Range("H2").Value = 0.08328125
$ws.Range("I2").Value = 0.18790625
$ws.Range("J2").Value = 0.3901875
$ws.Range("K2").Value = 0.3541875
$ws.Range("L2").Value = 0.60728125
$ws.Range("M2").Value = 0.11971875
$ws.Range("B3").Value = 0.2374375
$ws.Range("D3").Value = 0.8405625
$ws.Range("E3").Value = 0.77609375
$ws.Range("F3").Value = 0.8593125
$ws.Range("G3").Value = 0.8659375
$ws.Range("H3").Value = 0.8306249999999999
$ws.Range("I3").Value = 0.85615625
$ws.Range("J3").Value = 0.70665625
$ws.Range("K3").Value = 0.6058750000000001
$ws.Range("L3").Value = 0.5083124999999999
$ws.Range("M3").Value = 0.85453125
$ws.Range("B4").Value = 0.09193750000000001
$ws.Range("C4").Value = 0.8405625
$ws.Range("E4").Value = 0.98784375
$ws.Range("F4").Value = 0.97440625
$ws.Range("G4").Value = 0.98371875
$ws.Range("I4").Value = 0.95696875
$ws.Range("J4").Value = 0.7376875000000001
$ws.Range("K4").Value = 0.6464375
$ws.Range("L4").Value = 0.47153125
$ws.Range("M4").Value = 0.98825
$ws.Range("B5").Value = 0.04259375
$ws.Range("C5").Value = 0.77609375
$ws.Range("D5").Value = 0.98784375
$ws.Range("F5").Value = 0.7125312499999999
$ws.Range("G5").Value = 0.8783125000000001
$ws.Range("H5").Value = 0.99075
$ws.Range("I5").Value = 0.781125
$ws.Range("J5").Value = 0.56659375
$ws.Range("K5").Value = 0.6579375
$ws.Range("L5").Value = 0.3648125
$ws.Range("M5").Value = 0.9174375
$ws.Range("B6").Value = 0.19378125
$ws.Range("C6").Value = 0.8593125
$ws.Range("D6").Value = 0.97440625
$ws.Range("E6").Value = 0.7125312499999999
$ws.Range("G6").Value = 0.9931875
$ws.Range("H6").Value = 0.9488124999999999
$ws.Range("I6").Value = 0.99803125
$ws.Range("J6").Value = 0.8558125
$ws.Range("K6").Value = 0.61525
$ws.Range("L6").Value = 0.56534375
$ws.Range("M6").Value = 0.9660625
$ws.Range("B7").Value = 0.14053125
$ws.Range("C7").Value = 0.8659375
$ws.Range("D7").Value = 0.98371875
$ws.Range("E7").Value = 0.8783125000000001
$ws.Range("F7").Value = 0.9931875
$ws.Range("H7").Value = 0.96553125
$ws.Range("I7").Value = 0.9859375
$ws.Range("J7").Value = 0.7826875
$ws.Range("K7").Value = 0.64025
$ws.Range("L7").Value = 0.5070625
$ws.Range("M7").Value = 0.98796875
$ws.Range("B8").Value = 0.08328125
$ws.Range("C8").Value = 0.8306249999999999
$ws.Range("E8").Value = 0.99075
$ws.Range("F8").Value = 0.9488124999999999
$ws.Range("G8").Value = 0.96553125
$ws.Range("I8").Value = 0.93378125
$ws.Range("J8").Value = 0.71371875
$ws.Range("K8").Value = 0.6483125
$ws.Range("L8").Value = 0.454125
$ws.Range("M8").Value = 0.97915625
$ws.Range("B9").Value = 0.18790625
$ws.Range("C9").Value = 0.85615625
$ws.Range("D9").Value = 0.95696875
$ws.Range("E9").Value = 0.781125
$ws.Range("F9").Value = 0.99803125
$ws.Range("G9").Value = 0.9859375
$ws.Range("H9").Value = 0.93378125
$ws.Range("J9").Value = 0.83190625
$ws.Range("K9").Value = 0.6234062500000001
$ws.Range("L9").Value = 0.553375
$ws.Range("M9").Value = 0.95490625
$ws.Range("B10").Value = 0.3901875
$ws.Range("C10").Value = 0.70665625
$ws.Range("D10").Value = 0.7376875000000001
$ws.Range("E10").Value = 0.56659375
$ws.Range("F10").Value = 0.8558125
$ws.Range("G10").Value = 0.7826875
$ws.Range("H10").Value = 0.71371875
$ws.Range("I10").Value = 0.83190625
$ws.Range("K10").Value = 0.561875
$ws.Range("L10").Value = 0.58790625
$ws.Range("M10").Value = 0.73934375
$ws.Range("B11").Value = 0.3541875
$ws.Range("C11").Value = 0.6058750000000001
$ws.Range("D11").Value = 0.6464375
$ws.Range("E11").Value = 0.6579375
$ws.Range("F11").Value = 0.61525
$ws.Range("G11").Value = 0.64025
$ws.Range("H11").Value = 0.6483125
$ws.Range("I11").Value = 0.6234062500000001
$ws.Range("J11").Value = 0.561875
$ws.Range("L11").Value = 0.48721875
$ws.Range("M11").Value = 0.64396875
$ws.Range("B12").Value = 0.60728125
$ws.Range("C12").Value = 0.5083124999999999
$ws.Range("D12").Value = 0.47153125
$ws.Range("E12").Value = 0.3648125
$ws.Range("F12").Value = 0.56534375
$ws.Range("G12").Value = 0.5070625
$ws.Range("H12").Value = 0.454125
$ws.Range("I12").Value = 0.553375
$ws.Range("J12").Value = 0.58790625
$ws.Range("K12").Value = 0.48721875
$ws.Range("M12").Value = 0.48128125
$ws.Range("B13").Value = 0.11971875
$ws.Range("C13").Value = 0.85453125
$ws.Range("D13").Value = 0.98825
$ws.Range("E13").Value = 0.9174375
$ws.Range("F13").Value = 0.9660625
$ws.Range("G13").Value = 0.98796875
$ws.Range("H13").Value = 0.97915625
$ws.Range("I13").Value = 0.95490625
$ws.Range("J13").Value = 0.73934375
$ws.Range("K13").Value = 0.64396875
$ws.Range("L13").Value = 0.48128125

# Sheet 3: Pr(x<0)
$ws = $wb.Worksheets.Item(3)
$ws.Range("C2").Value = 0.7625625
$ws.Range("D2").Value = 0.9080625
$ws.Range("E2").Value = 0.95740625
$ws.Range("F2").Value = 0.80621875
$ws.Range("G2").Value = 0.85946875
$ws.Range("H2").Value = 0.91671875
$ws.Range("I2").Value = 0.8120937500000001
$ws.Range("J2").Value = 0.6098125
$ws.Range("K2").Value = 0.6458125
$ws.Range("L2").Value = 0.39271875
$ws.Range("M2").Value = 0.88028125
$ws.Range("B3").Value = 0.7625625
$ws.Range("D3").Value = 0.1594375
$ws.Range("E3").Value = 0.22390625
$ws.Range("F3").Value = 0.1406875
$ws.Range("G3").Value = 0.1340625
$ws.Range("H3").Value = 0.1693750000000001
$ws.Range("I3").Value = 0.14384375
$ws.Range("J3").Value = 0.29334375
$ws.Range("K3").Value = 0.3941249999999999
$ws.Range("L3").Value = 0.4916875000000001
$ws.Range("M3").Value = 0.14546875
$ws.Range("B4").Value = 0.9080625
$ws.Range("C4").Value = 0.1594375
$ws.Range("E4").Value = 0.01215624999999998
$ws.Range("F4").Value = 0.02559374999999997
$ws.Range("G4").Value = 0.01628125000000002
$ws.Range("I4").Value = 0.04303124999999997
$ws.Range("J4").Value = 0.2623124999999999
$ws.Range("K4").Value = 0.3535625
$ws.Range("L4").Value = 0.52846875
$ws.Range("M4").Value = 0.01175000000000004
$ws.Range("B5").Value = 0.95740625
$ws.Range("C5").Value = 0.22390625
$ws.Range("D5").Value = 0.01215624999999998
$ws.Range("F5").Value = 0.2874687500000001
$ws.Range("G5").Value = 0.1216874999999999
$ws.Range("H5").Value = 0.00924999999999998
$ws.Range("I5").Value = 0.218875
$ws.Range("J5").Value = 0.43340625
$ws.Range("K5").Value = 0.3420625
$ws.Range("L5").Value = 0.6351875
$ws.Range("M5").Value = 0.08256249999999998
$ws.Range("B6").Value = 0.80621875
$ws.Range("C6").Value = 0.1406875
$ws.Range("D6").Value = 0.02559374999999997
$ws.Range("E6").Value = 0.2874687500000001
$ws.Range("G6").Value = 0.006812499999999999
$ws.Range("H6").Value = 0.05118750000000005
$ws.Range("I6").Value = 0.001968749999999964
$ws.Range("J6").Value = 0.1441875
$ws.Range("K6").Value = 0.38475
$ws.Range("L6").Value = 0.43465625
$ws.Range("M6").Value = 0.03393749999999995
$ws.Range("B7").Value = 0.85946875
$ws.Range("C7").Value = 0.1340625
$ws.Range("D7").Value = 0.01628125000000002
$ws.Range("E7").Value = 0.1216874999999999
$ws.Range("F7").Value = 0.006812499999999999
$ws.Range("H7").Value = 0.03446875000000005
$ws.Range("I7").Value = 0.01406249999999998
$ws.Range("J7").Value = 0.2173125
$ws.Range("K7").Value = 0.35975
$ws.Range("L7").Value = 0.4929375
$ws.Range("M7").Value = 0.01203125000000005
$ws.Range("B8").Value = 0.91671875
$ws.Range("C8").Value = 0.1693750000000001
$ws.Range("E8").Value = 0.00924999999999998
$ws.Range("F8").Value = 0.05118750000000005
$ws.Range("G8").Value = 0.03446875000000005
$ws.Range("I8").Value = 0.06621874999999999
$ws.Range("J8").Value = 0.28628125
$ws.Range("K8").Value = 0.3516875
$ws.Range("L8").Value = 0.545875
$ws.Range("M8").Value = 0.02084375000000005
$ws.Range("B9").Value = 0.8120937500000001
$ws.Range("C9").Value = 0.14384375
$ws.Range("D9").Value = 0.04303124999999997
$ws.Range("E9").Value = 0.218875
$ws.Range("F9").Value = 0.001968749999999964
$ws.Range("G9").Value = 0.01406249999999998
$ws.Range("H9").Value = 0.06621874999999999
$ws.Range("J9").Value = 0.16809375
$ws.Range("K9").Value = 0.3765937499999999
$ws.Range("L9").Value = 0.446625
$ws.Range("M9").Value = 0.04509375000000004
$ws.Range("B10").Value = 0.6098125
$ws.Range("C10").Value = 0.29334375
$ws.Range("D10").Value = 0.2623124999999999
$ws.Range("E10").Value = 0.43340625
$ws.Range("F10").Value = 0.1441875
$ws.Range("G10").Value = 0.2173125
$ws.Range("H10").Value = 0.28628125
$ws.Range("I10").Value = 0.16809375
$ws.Range("K10").Value = 0.438125
$ws.Range("L10").Value = 0.41209375
$ws.Range("M10").Value = 0.26065625
$ws.Range("B11").Value = 0.6458125
$ws.Range("C11").Value = 0.3941249999999999
$ws.Range("D11").Value = 0.3535625
$ws.Range("E11").Value = 0.3420625
$ws.Range("F11").Value = 0.38475
$ws.Range("G11").Value = 0.35975
$ws.Range("H11").Value = 0.3516875
$ws.Range("I11").Value = 0.3765937499999999
$ws.Range("J11").Value = 0.438125
$ws.Range("L11").Value = 0.51278125
$ws.Range("M11").Value = 0.35603125
$ws.Range("B12").Value = 0.39271875
$ws.Range("C12").Value = 0.4916875000000001
$ws.Range("D12").Value = 0.52846875
$ws.Range("E12").Value = 0.6351875
$ws.Range("F12").Value = 0.43465625
$ws.Range("G12").Value = 0.4929375
$ws.Range("H12").Value = 0.545875
$ws.Range("I12").Value = 0.446625
$ws.Range("J12").Value = 0.41209375
$ws.Range("K12").Value = 0.51278125
$ws.Range("M12").Value = 0.51871875
$ws.Range("B13").Value = 0.88028125
$ws.Range("C13").Value = 0.14546875
$ws.Range("D13").Value = 0.01175000000000004
$ws.Range("E13").Value = 0.08256249999999998
$ws.Range("F13").Value = 0.03393749999999995
$ws.Range("G13").Value = 0.01203125000000005
$ws.Range("H13").Value = 0.02084375000000005
$ws.Range("I13").Value = 0.04509375000000004
$ws.Range("J13").Value = 0.26065625
$ws.Range("K13").Value = 0.35603125
$ws.Range("L13").Value = 0.51871875
